$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Column A looks like an ISO date ("2025-09-21"); Excel would normally parse
# that into a date serial on assignment. Prefix with an apostrophe so it is
# kept as literal text (as every other row in this sheet already is), then
# reset the cell style back to Normal so we don't leave a stray
# quote-prefix style behind.
$ws.Cells.Item($row, 1).Value = "'2025-09-21"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "21:20:53"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,777.8410"
